$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 245 (open/high/low/close corrected) ---
$ws.Range("C245").Value2 = 1568092000000
$ws.Range("D245").Value2 = 1568092000000
$ws.Range("E245").Value2 = 1568092000000
$ws.Range("F245").Value2 = 1568092000000

# --- Append three new monthly rows (253-255), matching the style of the last existing row ---
$ws.Range("A252").Copy()
$ws.Range("A253:A255").PasteSpecial(-4122)  # xlPasteFormats

# Row 253: 2023-01-01
$ws.Range("A253").Value2 = 44927.45833333334
$ws.Range("B253").Value2 = "ECONOMICS:AEM2"
$ws.Range("C253").Value2 = 1719863000000
$ws.Range("D253").Value2 = 1719863000000
$ws.Range("E253").Value2 = 1719863000000
$ws.Range("F253").Value2 = 1719863000000
$ws.Range("G253").Value2 = 0

# Row 254: 2023-02-01
$ws.Range("A254").Value2 = 44958.45833333334
$ws.Range("B254").Value2 = "ECONOMICS:AEM2"
$ws.Range("C254").Value2 = 1749400000000
$ws.Range("D254").Value2 = 1749400000000
$ws.Range("E254").Value2 = 1749400000000
$ws.Range("F254").Value2 = 1749400000000
$ws.Range("G254").Value2 = 0

# Row 255: 2023-03-01
$ws.Range("A255").Value2 = 44986.45833333334
$ws.Range("B255").Value2 = "ECONOMICS:AEM2"
$ws.Range("C255").Value2 = 1788400000000
$ws.Range("D255").Value2 = 1788400000000
$ws.Range("E255").Value2 = 1788400000000
$ws.Range("F255").Value2 = 1788400000000
$ws.Range("G255").Value2 = 0
